# Auto-generated Excel COM-interop script to apply the Golem_Profits diff
# Updates columns H-N (currentAveragePrice .. LeveProfitHQ) for specific rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (G19 context = 7015)
$ws.Range("H19").Value = 608.2353000000001
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 676
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = 676
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = -1026

# Row 28 (G28 context = 27772)
$ws.Range("H28").Value = 1573
$ws.Range("I28").Value = 1573
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1573
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -1088

# Row 43 (G43 context = 5472)
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").Value = 5000
$ws.Range("N43").Value = -5138

# Row 92 (G92 context = 19901)
$ws.Range("H92").Value = 316.42105
$ws.Range("I92").Value = 302.2143
$ws.Range("J92").Value = 356.2
$ws.Range("K92").Value = 302.2143
$ws.Range("L92").Value = 356.2
$ws.Range("M92").Value = 945.7857
$ws.Range("N92").Value = -2852.2

# Row 103 (G103 context = 19909)
$ws.Range("H103").Value = 500
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -914

# Row 107 (G107 context = 27766)
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()

# Row 118 (G118 context = 27958)
$ws.Range("H118").Value = 499.5
$ws.Range("I118").Value = 499.5
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1498.5
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 158.5


$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G2 context = 27713)
$ws.Range("H2").Value = 537.75
$ws.Range("I2").Value = 537.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 537.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -424.75

# Row 5 (G5 context = 5091)
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 29.142857
$ws.Range("J5").Value = 42
$ws.Range("K5").Value = 29.142857
$ws.Range("L5").Value = 42
$ws.Range("M5").Value = 82.85714300000001
$ws.Range("N5").Value = -266

# Row 55 (G55 context = 2830)
$ws.Range("H55").Value = 25500
$ws.Range("I55").Value = 25500
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 25500
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -25185

# Row 56 (G56 context = 2504)
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("N56").Value = 0

# Row 88 (G88 context = 12530)
$ws.Range("H88").Value = 3534
$ws.Range("I88").Value = 700
$ws.Range("J88").Value = 3938.8572
$ws.Range("K88").Value = 700
$ws.Range("L88").Value = 3938.8572
$ws.Range("M88").Value = -294
$ws.Range("N88").Value = -4750.8572

# Row 91 (G91 context = 12530)
$ws.Range("H91").Value = 3534
$ws.Range("I91").Value = 700
$ws.Range("J91").Value = 3938.8572
$ws.Range("K91").Value = 700
$ws.Range("L91").Value = 3938.8572
$ws.Range("M91").Value = 704
$ws.Range("N91").Value = -6746.8572

# Row 110 (G110 context = 27708)
$ws.Range("H110").Value = 999.5
$ws.Range("I110").Value = 999.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 999.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1045.5

# Row 116 (G116 context = 27713)
$ws.Range("H116").Value = 537.75
$ws.Range("I116").Value = 537.75
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 537.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1756.25

# Row 124 (G124 context = 34252)
$ws.Range("H124").Value = 34688.8
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 34688.8
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 34688.8
$ws.Range("N124").Value = -44508.8

# Row 125 (G125 context = 34251)
$ws.Range("H125").Value = 62857.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 62857.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 62857.5
$ws.Range("N125").Value = -72697.5

# Row 132 (G132 context = 43997)
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 10500
$ws.Range("N132").Value = -15560


$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G3 context = 27713)
$ws.Range("H3").Value = 537.75
$ws.Range("I3").Value = 537.75
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 537.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -423.75

# Row 4 (G4 context = 5091)
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 29.142857
$ws.Range("J4").Value = 42
$ws.Range("K4").Value = 29.142857
$ws.Range("L4").Value = 42
$ws.Range("M4").Value = 85.85714300000001
$ws.Range("N4").Value = -272

# Row 20 (G20 context = 14149)
$ws.Range("H20").Value = 1900
$ws.Range("I20").Value = 1900
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1900
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1653

# Row 36 (G36 context = 2320)
$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 34


$ws = $wb.Worksheets.Item("CRP")
# Row 5 (G5 context = 1893)
$ws.Range("H5").Value = 3025
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 700
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -588
$ws.Range("N5").Value = -10224

# Row 16 (G16 context = 27691)
$ws.Range("H16").Value = 881.5
$ws.Range("I16").Value = 657.8
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 657.8
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -370.8
$ws.Range("N16").Value = -2574

# Row 22 (G22 context = 5367)
$ws.Range("H22").Value = 782.8570999999999
$ws.Range("I22").Value = 782.8570999999999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 782.8570999999999
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -432.8570999999999

# Row 58 (G58 context = 44021)
$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2000
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -1797

# Row 59 (G59 context = 1942)
$ws.Range("H59").Value = 46276
$ws.Range("I59").Value = 104
$ws.Range("J59").Value = 61666.668
$ws.Range("K59").Value = 104
$ws.Range("L59").Value = 61666.668
$ws.Range("M59").Value = 1041
$ws.Range("N59").Value = -63956.668

# Row 60 (G60 context = 1937)
$ws.Range("H60").Value = 28644.285
$ws.Range("I60").Value = 12060.667
$ws.Range("J60").Value = 41082
$ws.Range("K60").Value = 12060.667
$ws.Range("L60").Value = 41082
$ws.Range("M60").Value = -11549.667
$ws.Range("N60").Value = -42104

# Row 103 (G103 context = 19558)
$ws.Range("H103").Value = 48500
$ws.Range("I103").Value = 48500
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 48500
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -47328

# Row 105 (G105 context = 19928)
$ws.Range("H105").Value = 615.8
$ws.Range("I105").Value = 519.75
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 519.75
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1227.25
$ws.Range("N105").Value = -4494

# Row 113 (G113 context = 27691)
$ws.Range("H113").Value = 881.5
$ws.Range("I113").Value = 657.8
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 657.8
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1512.2
$ws.Range("N113").Value = -6340

# Row 132 (G132 context = 44019)
$ws.Range("H132").Value = 1210
$ws.Range("I132").Value = 1210
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3630
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1100

# Row 136 (G136 context = 44021)
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -3450


$ws = $wb.Worksheets.Item("CUL")
# Row 2 (G2 context = 4847)
$ws.Range("H2").Value = 13
$ws.Range("I2").Value = 7.7647057
$ws.Range("J2").Value = 22.88889
$ws.Range("K2").Value = 46.5882342
$ws.Range("L2").Value = 137.33334
$ws.Range("M2").Value = 66.4117658
$ws.Range("N2").Value = -363.33334


$ws = $wb.Worksheets.Item("GSM")
# Row 113 (G113 context = 27710)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# Row 122 (G122 context = 36182)
$ws.Range("H122").Value = 3854.25
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 6420
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 19260
$ws.Range("M122").Value = -6547
$ws.Range("N122").Value = -24160


$ws = $wb.Worksheets.Item("LTW")
# Row 16 (G16 context = 5289)
$ws.Range("H16").Value = 5176.5
$ws.Range("I16").Value = 10000
$ws.Range("J16").Value = 353
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 353
$ws.Range("M16").Value = -9830
$ws.Range("N16").Value = -693

# Row 40 (G40 context = 36248)
$ws.Range("H40").Value = 651000.1
$ws.Range("I40").Value = 51997.5
$ws.Range("J40").Value = 850667.7
$ws.Range("K40").Value = 51997.5
$ws.Range("L40").Value = 850667.7
$ws.Range("M40").Value = -51861.5
$ws.Range("N40").Value = -850939.7

# Row 100 (G100 context = 19995)
$ws.Range("H100").Value = 2999
$ws.Range("I100").Value = 2999
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2999
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2458

# Row 104 (G104 context = 18675)
$ws.Range("H104").Value = 26252.834
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 26252.834
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 26252.834
$ws.Range("N104").Value = -33240.834


$ws = $wb.Worksheets.Item("WVR")
# Row 81 (G81 context = 12596)
$ws.Range("H81").Value = 1650
$ws.Range("I81").Value = 1650
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3300
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -2239

# Row 84 (G84 context = 12596)
$ws.Range("H84").Value = 1650
$ws.Range("I84").Value = 1650
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 16500
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -11196

# Row 100 (G100 context = 19981)
$ws.Range("H100").Value = 5002
$ws.Range("I100").Value = 5002
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 10004
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -9463

# Row 113 (G113 context = 27752)
$ws.Range("H113").Value = 8462.429
$ws.Range("I113").Value = 560.5
$ws.Range("J113").Value = 18998.334
$ws.Range("K113").Value = 1681.5
$ws.Range("L113").Value = 56995.00199999999
$ws.Range("M113").Value = 488.5
$ws.Range("N113").Value = -61335.00199999999

